# Driver timesheet report added
# Populate the "Лист1" (sheet 1) report header cells and size the columns
# that make up the new driver-timesheet layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# -- Report header text ------------------------------------------------
$ws.Range("B2").Value = "Номер маршрута:"
$ws.Range("B4").Value = "Дата создания отчета: "

# -- Column widths for the new layout -----------------------------------
$ws.Columns.Item(2).ColumnWidth = 17            # B: route number label
$ws.Columns.Item(3).ColumnWidth = 20            # C: value column
$ws.Range($ws.Columns.Item(4), $ws.Columns.Item(9)).ColumnWidth = 8   # D:I
$ws.Range($ws.Columns.Item(35), $ws.Columns.Item(36)).ColumnWidth = 12 # AI:AJ
$ws.Columns.Item(37).ColumnWidth = 17            # AK

# -- Selection moves to C11 as left by the report author ----------------
[void]$ws.Range("C11").Select()
